$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template row (carries number format / border / font style for column A)
$src = $ws.Range("A233:D233")

# New data rows to append: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @{ Row = 234; A = 44308; B = 1; C = 3; D = 131.3485113835376 },
    @{ Row = 235; A = 44309; B = 0; C = 3; D = 131.3485113835376 },
    @{ Row = 236; A = 44310; B = 0; C = 3; D = 131.3485113835376 },
    @{ Row = 237; A = 44311; B = 0; C = 3; D = 131.3485113835376 },
    @{ Row = 238; A = 44312; B = 0; C = 3; D = 131.3485113835376 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $dst = $ws.Range("A" + $rowNum + ":D" + $rowNum)
    $src.Copy($dst)

    $ws.Range("A" + $rowNum).Value = $r.A
    $ws.Range("B" + $rowNum).Value = $r.B
    $ws.Range("C" + $rowNum).Value = $r.C
    $ws.Range("D" + $rowNum).Value = $r.D
}
